$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 1).Value = 1
$ws.Cells.Item(2, 3).Value = 4
$ws.Cells.Item(2, 4).Value = 3
$ws.Cells.Item(2, 5).Value = 0
$ws.Cells.Item(2, 6).Value = 1
$ws.Cells.Item(2, 7).Value = 2
$ws.Cells.Item(2, 8).Value = 1
$ws.Cells.Item(2, 9).Value = 1
$ws.Cells.Item(2, 10).Value = 3
$ws.Cells.Item(2, 11).Value = 0
$ws.Cells.Item(2, 12).Value = 0
$ws.Cells.Item(2, 13).Value = 9

$ws.Cells.Item(3, 1).Value = 2
$ws.Cells.Item(3, 3).Value = 4
$ws.Cells.Item(3, 4).Value = 1
$ws.Cells.Item(3, 5).Value = 0
$ws.Cells.Item(3, 6).Value = 3
$ws.Cells.Item(3, 7).Value = 2
$ws.Cells.Item(3, 8).Value = 6
$ws.Cells.Item(3, 9).Value = -4
$ws.Cells.Item(3, 10).Value = 1
$ws.Cells.Item(3, 11).Value = 0
$ws.Cells.Item(3, 12).Value = 0
$ws.Cells.Item(3, 13).Value = 3

$ws.Cells.Item(4, 1).Value = 3
$ws.Cells.Item(4, 3).Value = 4
$ws.Cells.Item(4, 4).Value = 4
$ws.Cells.Item(4, 5).Value = 0
$ws.Cells.Item(4, 6).Value = 0
$ws.Cells.Item(4, 7).Value = 5
$ws.Cells.Item(4, 8).Value = 0
$ws.Cells.Item(4, 9).Value = 5
$ws.Cells.Item(4, 10).Value = 2
$ws.Cells.Item(4, 11).Value = 1
$ws.Cells.Item(4, 12).Value = 0
$ws.Cells.Item(4, 13).Value = 13

$ws.Cells.Item(5, 1).Value = 4
$ws.Cells.Item(5, 3).Value = 4
$ws.Cells.Item(5, 4).Value = 0
$ws.Cells.Item(5, 5).Value = 1
$ws.Cells.Item(5, 6).Value = 3
$ws.Cells.Item(5, 7).Value = 0
$ws.Cells.Item(5, 8).Value = 3
$ws.Cells.Item(5, 9).Value = -3
$ws.Cells.Item(5, 10).Value = 0
$ws.Cells.Item(5, 11).Value = 0
$ws.Cells.Item(5, 12).Value = 0
$ws.Cells.Item(5, 13).Value = 1

$ws.Cells.Item(6, 1).Value = 5
$ws.Cells.Item(6, 3).Value = 4
$ws.Cells.Item(6, 4).Value = 3
$ws.Cells.Item(6, 5).Value = 0
$ws.Cells.Item(6, 6).Value = 1
$ws.Cells.Item(6, 7).Value = 7
$ws.Cells.Item(6, 8).Value = 2
$ws.Cells.Item(6, 9).Value = 5
$ws.Cells.Item(6, 10).Value = 1
$ws.Cells.Item(6, 11).Value = 0
$ws.Cells.Item(6, 12).Value = 2
$ws.Cells.Item(6, 13).Value = 13

$ws.Cells.Item(7, 1).Value = 6
$ws.Cells.Item(7, 3).Value = 4
$ws.Cells.Item(7, 4).Value = 1
$ws.Cells.Item(7, 5).Value = 1
$ws.Cells.Item(7, 6).Value = 2
$ws.Cells.Item(7, 7).Value = 1
$ws.Cells.Item(7, 8).Value = 2
$ws.Cells.Item(7, 9).Value = -1
$ws.Cells.Item(7, 10).Value = 1
$ws.Cells.Item(7, 11).Value = 0
$ws.Cells.Item(7, 12).Value = 0
$ws.Cells.Item(7, 13).Value = 4

$ws.Cells.Item(8, 1).Value = 7
$ws.Cells.Item(8, 3).Value = 4
$ws.Cells.Item(8, 4).Value = 1
$ws.Cells.Item(8, 5).Value = 0
$ws.Cells.Item(8, 6).Value = 3
$ws.Cells.Item(8, 7).Value = 1
$ws.Cells.Item(8, 8).Value = 5
$ws.Cells.Item(8, 9).Value = -4
$ws.Cells.Item(8, 10).Value = 1
$ws.Cells.Item(8, 11).Value = 0
$ws.Cells.Item(8, 12).Value = 0
$ws.Cells.Item(8, 13).Value = 3

$ws.Cells.Item(9, 1).Value = 8
$ws.Cells.Item(9, 3).Value = 4
$ws.Cells.Item(9, 4).Value = 2
$ws.Cells.Item(9, 5).Value = 0
$ws.Cells.Item(9, 6).Value = 2
$ws.Cells.Item(9, 7).Value = 2
$ws.Cells.Item(9, 8).Value = 2
$ws.Cells.Item(9, 9).Value = 0
$ws.Cells.Item(9, 10).Value = 2
$ws.Cells.Item(9, 11).Value = 0
$ws.Cells.Item(9, 12).Value = 0
$ws.Cells.Item(9, 13).Value = 6

$lo = $ws.ListObjects.Item(1)
$lo.Name = "Table_1"

